$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to write a "Conta" (account number) value while preserving leading
# zeros: these look like numbers but must stay as text, so force the cell's
# number format to Text before assigning the string.
function Set-ContaText($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# 1) Row 5: GABRIELA/10000 -> LEILA/9000 (simple in-place edit, no shift)
Set-ContaText $ws.Cells.Item(5, 1) "004208447"
$ws.Cells.Item(5, 2).Value = "LEILA"
$ws.Cells.Item(5, 3).Value = 9000

# 2) Insert a brand-new row at position 9 (pushes AHMAD and everything below
#    down by one) and populate it with DIOGO/3000.
$ws.Rows.Item(9).Insert()
Set-ContaText $ws.Cells.Item(9, 1) "004550415"
$ws.Cells.Item(9, 2).Value = "DIOGO"
$ws.Cells.Item(9, 3).Value = 3000

# 3) Row 13 (originally row 12, ELENE/2163.87, now shifted down by the
#    insert above): ELENE/2163.87 -> CESAR/2000
Set-ContaText $ws.Cells.Item(13, 1) "004207278"
$ws.Cells.Item(13, 2).Value = "CESAR"
$ws.Cells.Item(13, 3).Value = 2000

# 4) Insert another brand-new row right after, at position 14, and populate
#    it with JOSE/1500.
$ws.Rows.Item(14).Insert()
Set-ContaText $ws.Cells.Item(14, 1) "004480134"
$ws.Cells.Item(14, 2).Value = "JOSE"
$ws.Cells.Item(14, 3).Value = 1500
